# Adds two new attendance columns ("Telat" and "Cuti") to the "Data Karyawan"
# header table, shifting "Tidak Hadir" one column to the right, and refreshes
# the header formatting (center/center alignment, borderless title row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title row (row 1): drop the bottom border, keep everything else ---
$a1 = $ws.Cells.Item(1, 1)
$a1.Borders.LineStyle = 0
$a1.Copy() | Out-Null
$ws.Range("B1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- A2:I2 (plain header cells): vertical alignment top -> center ---
$a2 = $ws.Cells.Item(2, 1)
$a2.HorizontalAlignment = -4108   # xlCenter
$a2.VerticalAlignment = -4108     # xlCenter
$a2.Copy() | Out-Null
$ws.Range("B2:I2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- "Orange" header style reference (fill #FFFF3300) -> center/center ---
$j2 = $ws.Cells.Item(2, 10)   # Tanggal Mulai Kerja
$j2.HorizontalAlignment = -4108
$j2.VerticalAlignment = -4108
$j2.Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null    # Waktu Awal Kerja
$excel.CutCopyMode = 0

# --- "Blue" header style reference (fill #FF00B0F0) -> center/center ---
$k2 = $ws.Cells.Item(2, 11)   # Tanggal Akhir Kerja
$k2.HorizontalAlignment = -4108
$k2.VerticalAlignment = -4108
$k2.Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null    # Waktu Selesai kerja
$excel.CutCopyMode = 0

# --- Shift "Tidak Hadir" from O2 to P2 to make room for "Telat" ---
$ws.Cells.Item(2, 15).Cut($ws.Cells.Item(2, 16)) | Out-Null

# --- New header text ---
$ws.Cells.Item(2, 15).Value = "Telat"   # O2 (new)
$ws.Cells.Item(2, 17).Value = "Cuti"    # Q2 (new)

# --- Re-apply the alternating orange/blue look across N2:Q2 ---
$j2.Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null   # Hadir       -> orange
$j2.Copy() | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null   # Tidak Hadir -> orange
$k2.Copy() | Out-Null
$ws.Range("O2").PasteSpecial(-4122) | Out-Null   # Telat       -> blue
$k2.Copy() | Out-Null
$ws.Range("Q2").PasteSpecial(-4122) | Out-Null   # Cuti        -> blue
$excel.CutCopyMode = 0

# Make sure the textual values are exactly right (PasteSpecial only touches formats)
$ws.Cells.Item(2, 14).Value2 = "Hadir"
$ws.Cells.Item(2, 15).Value2 = "Telat"
$ws.Cells.Item(2, 16).Value2 = "Tidak Hadir"
$ws.Cells.Item(2, 17).Value2 = "Cuti"

# --- New column (P) width, manually resized by the author ---
$ws.Columns.Item(16).ColumnWidth = 13.1796875

# --- Selection bookkeeping, matches the saved view state ---
$ws.Range("G5").Select() | Out-Null
